$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "current page (bookmark)" for "The Spirit of Kaizen" (row 14) from 27 to 30
$ws.Range("C14").Value = 30

# Update the selected cell to reflect where the user left off editing
$ws.Range("C17").Select()
